# Add the new book entry "This is Going to Hurt" by Adam Kay to the
# "Completed" reading list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# New row is the row right after the last used row.
$newRow = 42

$ws.Cells.Item($newRow, 1).Value = "This is Going to Hurt"
$ws.Cells.Item($newRow, 2).Value = "Adam Kay"

# Start Date / Finish Date - both serial date 43912 (2020-03-22).
# Copy the formatting (date number format) from the row above first,
# then set the actual values so no new custom number-format style gets
# introduced.
$ws.Cells.Item($newRow - 1, 3).Copy($ws.Cells.Item($newRow, 3))
$ws.Cells.Item($newRow - 1, 4).Copy($ws.Cells.Item($newRow, 4))

$bookDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(43912)
$ws.Cells.Item($newRow, 3).Value = $bookDate
$ws.Cells.Item($newRow, 4).Value = $bookDate

$ws.Cells.Item($newRow, 5).Value = "medicine;doctor;nhs;burnout;science"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "5 Hours 53 Mins"

# Update the view to match what Excel recorded after the edit: the
# window was scrolled down one row and the next empty cell (B43) was
# selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("B43").Select()
